$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text representation
# (values like "139.00", "0.06150", "26.004.56" must not be
# auto-converted to numbers, which would drop formatting/precision)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.004.56"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.751.43"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "235.67"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5201"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("D9").Value = "39.64"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "0.06150"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "1.755.76"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "15.51"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "0.6453"
$ws.Range("E14").Value = "  +5.46%  "
$ws.Range("D15").Value = "4.525"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "77.58"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "26.009.11"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "0.000006633"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "1.977.13"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").Value = "8.651"
$ws.Range("E24").Value = "  +5.33%  "
$ws.Range("D25").Value = "5.164"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "139.00"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "1.507"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("D28").Value = "1.848"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "15.10"
$ws.Range("D30").Value = "103.08"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "0.08319"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "3.657"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "3.440"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "0.04445"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").Value = "2.609"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "0.9883"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "0.6102"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "0.01586"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").Value = "1.948"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "0.9981"
$ws.Range("D42").Value = "100.77"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "0.3873"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "0.7352"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").Value = "5.036"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").Value = "0.05477"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "6.352"
$ws.Range("E47").Value = "  +7.15%  "
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").Value = "52.93"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").Value = "30.06"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "7.558"
$ws.Range("E51").Value = "  +1.46%  "
